$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get permuted across rows: D, M, N, O, P, S
$cols = @("D", "M", "N", "O", "P", "S")

# Mapping: new row index -> old row index (data to be copied into new row comes from old row)
$mapping = @{
    2  = 11
    3  = 8
    4  = 9
    5  = 5
    6  = 4
    7  = 7
    8  = 6
    9  = 2
    10 = 3
    11 = 10
}

# Capture old values for rows 2..11 for the affected columns before any write
$oldValues = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $oldValues[$r] = $rowVals
}

# Apply new values according to the mapping
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $oldValues[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
